$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MuSCs -> Edn3 -> Ednrb -> ECs)
$ws.Range("G2").Value = 7.736660666666666
$ws.Range("H2").Value = 23.209982
$ws.Range("M2").Value = 52.271196
$ws.Range("N2").Value = 156.813588
$ws.Range("O2").Value = 0.6500553798777896
$ws.Range("P2").Value = 0.6500553798777895
$ws.Range("Q2").Value = 404.404506092824
$ws.Range("R2").Value = 3639.640554835416
$ws.Range("S2").Value = 0.6500553798777896
$ws.Range("T2").Value = 0.6500553798777895

# Row 3 (MuSCs -> Edn3 -> Ednrb -> FAPs)
$ws.Range("G3").Value = 7.736660666666666
$ws.Range("H3").Value = 23.209982
$ws.Range("O3").Value = 0.001694346062422021
$ws.Range("P3").Value = 0.00169434606242202
$ws.Range("Q3").Value = 1.054065859208667
$ws.Range("R3").Value = 9.486592732878
$ws.Range("S3").Value = 0.001694346062422021
$ws.Range("T3").Value = 0.00169434606242202

# Row 4 (MuSCs -> Edn3 -> Ednrb -> MuSCs)
$ws.Range("G4").Value = 7.736660666666666
$ws.Range("H4").Value = 23.209982
$ws.Range("M4").Value = 24.41792966666667
$ws.Range("N4").Value = 73.253789
$ws.Range("O4").Value = 0.3036664120961408
$ws.Range("P4").Value = 0.3036664120961408
$ws.Range("Q4").Value = 188.9132360135331
$ws.Range("R4").Value = 1700.219124121798
$ws.Range("S4").Value = 0.3036664120961408
$ws.Range("T4").Value = 0.3036664120961408

# Row 5 (MuSCs -> Edn3 -> Ednrb -> Resolving-Mac)
$ws.Range("G5").Value = 7.736660666666666
$ws.Range("H5").Value = 23.209982
$ws.Range("M5").Value = 3.585005
$ws.Range("N5").Value = 10.755015
$ws.Range("O5").Value = 0.04458386196364773
$ws.Range("P5").Value = 0.04458386196364771
$ws.Range("Q5").Value = 27.73596717330333
$ws.Range("R5").Value = 249.62370455973
$ws.Range("S5").Value = 0.04458386196364773
$ws.Range("T5").Value = 0.04458386196364771
